# "Fruta / hortaliza, semanal" — weekly refresh of the Espinaca price sheet.
#
# A new weekly price observation (fecha 2021-11-12) is inserted as row 22,
# pushing the previously existing rows 22:49 down to 23:50 (dimension grows
# from A1:R49 to A1:R50). The new row reuses the same commercial-unit /
# origin / classification values as the existing "$6.500-7.000, cuna 10
# kilos, Región Metropolitana" entries (the row that ends up at row 25),
# only the date differs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 22 — Excel shifts rows 22:49 down to 23:50,
# extending the sheet's used range to A1:R50.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22.
$ws.Cells.Item(22, 1).Value  = 11
$ws.Cells.Item(22, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(22, 3).Value  = "Bíobío"
$ws.Cells.Item(22, 4).Value  = 44512
$ws.Cells.Item(22, 5).Value  = 8
$ws.Cells.Item(22, 6).Value  = 100112012
$ws.Cells.Item(22, 7).Value  = "Espinaca"
$ws.Cells.Item(22, 8).Value  = "Sin especificar"
$ws.Cells.Item(22, 9).Value  = "Primera"
$ws.Cells.Item(22, 10).Value = 60
$ws.Cells.Item(22, 11).Value = 6500
$ws.Cells.Item(22, 12).Value = 7000
$ws.Cells.Item(22, 13).Value = 6750
$ws.Cells.Item(22, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(22, 15).Value = "Región Metropolitana"
$ws.Cells.Item(22, 16).Value = 675
$ws.Cells.Item(22, 17).Value = 10
$ws.Cells.Item(22, 18).Value = "Hortaliza"
